$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# D15: update value
$ws.Range("D15").Value = 300000000

# D16: update value
$ws.Range("D16").Value = -506013111.95999998

# D18: set formula (was a plain number)
$ws.Range("D18").Formula = "=SUM(D12:D17)"

# D19: was blank, now has a value
$ws.Range("D19").Value = -383100000

# D21: was blank, now has a formula
$ws.Range("D21").Formula = "=SUM(D18:D20)"

# D22: was blank, now has a value
$ws.Range("D22").Value = -32201025

# D23 already has formula =SUM(D21:D22); just let it recalc

# D25 already has formula =D23/D26; previously errored because D26 was text
# D26: was text " " (space), now a numeric value
$ws.Range("D26").Value = 1009991810

# D28, E28, C29 already have formulas; they will recalc once dependencies fixed

$excel.Calculate()
